$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = 641

$ws.Range("R3").Value = 7392
$ws.Range("S3").Value = 2672

$ws.Range("R4").Value = 1477
$ws.Range("S4").Value = 254

$ws.Range("R5").Value = 3502
$ws.Range("S5").Value = 131

$ws.Range("R6").Value = 4348
$ws.Range("S6").Value = 974

$ws.Range("R7").Value = 3191
$ws.Range("S7").Value = 138

$ws.Range("R8").Value = 3499
$ws.Range("S8").Value = 174

$ws.Range("R9").Value = 7645
$ws.Range("S9").Value = 2685

$ws.Range("R15").Select()
